# ---------------------------------------------------------------------------
# Edit: (1) change the table style on the "Plenary - complete the missing
#           gaps" slide (slide 16) from the custom green style to the
#           built-in style {B4B32D4F-6C6F-44C8-854B-D3C9E13AB54F};
#       (2) switch the deck's theme palette from the custom "Integral"
#           colours over to the standard "Office" palette (dk1/lt1/dk2/lt2/
#           accent1-6/hlink/folHlink), matching a Design > Theme change.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style -----------------------------------------------------
$s   = $p.Slides.Item(16)
$tbl = $s.Shapes.Item(3).Table
$tbl.ApplyStyle("{B4B32D4F-6C6F-44C8-854B-D3C9E13AB54F}")

# --- 2. Theme colours ----------------------------------------------------
# RGB() isn't available in this host, so colours are passed as the packed
# integer PowerPoint itself stores (r + g*256 + b*65536).
$scheme = $p.SlideMaster.ColorScheme
$scheme.Colors(1).RGB  = 0         # dk1      000000
$scheme.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$scheme.Colors(3).RGB  = 6968388   # dk2      44546A
$scheme.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$scheme.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$scheme.Colors(6).RGB  = 3243501   # accent2  ED7D31
$scheme.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$scheme.Colors(8).RGB  = 49407     # accent4  FFC000
$scheme.Colors(9).RGB  = 12874308  # accent5  4472C4
$scheme.Colors(10).RGB = 4697456   # accent6  70AD47
$scheme.Colors(11).RGB = 12673797  # hlink    0563C1
$scheme.Colors(12).RGB = 7491477   # folHlink 954F72
